$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# 1) Rows 1-3: single-value cells become "0M"
$t.Cell(1, 1).Range.Text = "0M"
$t.Cell(2, 1).Range.Text = "0M"
$t.Cell(3, 1).Range.Text = "0M"

# 2) Insert 10 new rows right after row 3 (i.e. before the row that is
#    currently row 4), each holding one of the new values, in order.
$newValues = @("2538", "0.00003", "0.00008", "0.00004", "0.00001", "0.00004", "0.00004", "0.00005", "0.11345", "100.0")

$anchorIndex = 4
foreach ($val in $newValues) {
    $newRow = $t.Rows.Add($t.Rows.Item($anchorIndex))
    $t.Cell($newRow.Index, 1).Range.Text = $val
    $anchorIndex = $newRow.Index + 1
}

# 3) The last three rows (originally tab-separated multi-run cells) collapse
#    down to a single value each.
$totalRows = $t.Rows.Count
$t.Cell($totalRows - 2, 1).Range.Text = "99.96"
$t.Cell($totalRows - 1, 1).Range.Text = "0.11"
$t.Cell($totalRows, 1).Range.Text = "270"
